# Applies the record re-sync update for rows 16-22 and 24-33 of the
# "Artfynd" sheet, matching the target OOXML diff exactly (cell-by-cell).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16
$ws.Range("A16").Value = 111814688
$ws.Range("B16").Value = 90087
$ws.Range("D16").Value = 'LC'
$ws.Range("E16").Value = 3298
$ws.Range("F16").Value = 'Trådticka'
$ws.Range("G16").Value = 'Climacocystis borealis'
$ws.Range("H16").Value = '(Fr.) Kotl. & Pouzar'
$ws.Range("P16").Value = 'åsele 1:1 (åsele 1:1), Ås lm'
$ws.Range("Q16").Value = 610011
$ws.Range("R16").Value = 7121476
$ws.Range("S16").Value = 1
$ws.Range("Z16").Value = '17:55'
$ws.Range("AB16").Value = '17:55'

# Row 17
$ws.Range("A17").Value = 111814591
$ws.Range("B17").Value = 77515
$ws.Range("D17").Value = 'NT'
$ws.Range("E17").Value = 6425
$ws.Range("F17").Value = 'Garnlav'
$ws.Range("G17").Value = 'Alectoria sarmentosa'
$ws.Range("H17").Value = '(Ach.) Ach.'
$ws.Range("Q17").Value = 610012
$ws.Range("R17").Value = 7121464
$ws.Range("Z17").Value = '17:50'
$ws.Range("AB17").Value = '17:50'

# Row 18
$ws.Range("A18").Value = 111814104
$ws.Range("B18").Value = 56398
$ws.Range("E18").Value = 100109
$ws.Range("F18").Value = 'Tretåig hackspett'
$ws.Range("G18").Value = 'Picoides tridactylus'
$ws.Range("H18").Value = '(Linnaeus, 1758)'
$ws.Range("Q18").Value = 610155
$ws.Range("R18").Value = 7121460
$ws.Range("Z18").Value = '17:23'
$ws.Range("AB18").Value = '17:23'

# Row 19
$ws.Range("A19").Value = 111815024
$ws.Range("B19").Value = 56414
$ws.Range("E19").Value = 100049
$ws.Range("F19").Value = 'Spillkråka'
$ws.Range("G19").Value = 'Dryocopus martius'
$ws.Range("H19").Value = '(Linnaeus, 1758)'
$ws.Range("Q19").Value = 609922
$ws.Range("R19").Value = 7121488
$ws.Range("Z19").Value = '18:12'
$ws.Range("AB19").Value = '18:12'

# Row 20
$ws.Range("A20").Value = 111814478
$ws.Range("B20").Value = 77515
$ws.Range("D20").Value = 'NT'
$ws.Range("E20").Value = 6425
$ws.Range("F20").Value = 'Garnlav'
$ws.Range("G20").Value = 'Alectoria sarmentosa'
$ws.Range("H20").Value = '(Ach.) Ach.'
$ws.Range("Q20").Value = 610155
$ws.Range("R20").Value = 7121461
$ws.Range("Z20").Value = '17:41'
$ws.Range("AB20").Value = '17:41'

# Row 21
$ws.Range("A21").Value = 111815269
$ws.Range("B21").Value = 90666
$ws.Range("D21").Value = 'LC'
$ws.Range("E21").Value = 4364
$ws.Range("F21").Value = 'Dropptaggsvamp'
$ws.Range("G21").Value = 'Hydnellum ferrugineum'
$ws.Range("H21").Value = '(Fr.:Fr.) P. Karst.'
$ws.Range("Q21").Value = 610054
$ws.Range("R21").Value = 7121273
$ws.Range("Z21").Value = '18:27'
$ws.Range("AB21").Value = '18:27'

# Row 22
$ws.Range("A22").Value = 111815114
$ws.Range("B22").Value = 90660
$ws.Range("E22").Value = 4362
$ws.Range("F22").Value = 'Blå taggsvamp'
$ws.Range("G22").Value = 'Hydnellum caeruleum'
$ws.Range("H22").Value = '(Hornem.) P.Karst.'
$ws.Range("P22").Value = 'åsele 1:1, Ås lm'
$ws.Range("Q22").Value = 610384
$ws.Range("R22").Value = 7121170
$ws.Range("S22").Value = 5
$ws.Range("Z22").Value = '18:19'
$ws.Range("AB22").Value = '18:19'

# Row 24
$ws.Range("B24").Value = 88623

# Row 25
$ws.Range("A25").Value = 112013700
$ws.Range("B25").Value = 77636
$ws.Range("Q25").Value = 610102
$ws.Range("R25").Value = 7121416
$ws.Range("Z25").Value = '19:35'
$ws.Range("AB25").Value = '19:35'

# Row 26
$ws.Range("A26").Value = 112013704
$ws.Range("B26").Value = 81371
$ws.Range("E26").Value = 1312
$ws.Range("F26").Value = 'Gammelgransskål'
$ws.Range("G26").Value = 'Pseudographis pinicola'
$ws.Range("H26").Value = '(Nyl.) Rehm'
$ws.Range("Q26").Value = 610094
$ws.Range("R26").Value = 7121455
$ws.Range("Z26").Value = '19:49'
$ws.Range("AB26").Value = '19:49'

# Row 27
$ws.Range("A27").Value = 112013697
$ws.Range("B27").Value = 89557
$ws.Range("E27").Value = 5432
$ws.Range("F27").Value = 'Granticka'
$ws.Range("G27").Value = 'Porodaedalea chrysoloma'
$ws.Range("H27").Value = '(Fr.) Fiasson & Niemelä'
$ws.Range("Q27").Value = 610102
$ws.Range("R27").Value = 7121413
$ws.Range("Z27").Value = '19:35'
$ws.Range("AB27").Value = '19:35'

# Row 28
$ws.Range("A28").Value = 112013690
$ws.Range("B28").Value = 88623
$ws.Range("E28").Value = 1962
$ws.Range("F28").Value = 'Vaddporing'
$ws.Range("G28").Value = 'Anomoporia kamtschatica'
$ws.Range("H28").Value = '(Parmasto) Bondartseva'
$ws.Range("Q28").Value = 610052
$ws.Range("R28").Value = 7121425
$ws.Range("Z28").Value = '19:43'
$ws.Range("AB28").Value = '19:43'

# Row 29
$ws.Range("A29").Value = 112013703
$ws.Range("Q29").Value = 610144
$ws.Range("R29").Value = 7121461
$ws.Range("Z29").Value = '19:28'
$ws.Range("AB29").Value = '19:28'

# Row 30
$ws.Range("A30").Value = 112013696
$ws.Range("B30").Value = 87095
$ws.Range("E30").Value = 4962
$ws.Range("F30").Value = 'Mjölsvärting'
$ws.Range("G30").Value = 'Lyophyllum semitale'
$ws.Range("H30").Value = '(Fr. : Fr.) Kühner'
$ws.Range("Q30").Value = 610070
$ws.Range("R30").Value = 7121402
$ws.Range("Z30").Value = '19:40'
$ws.Range("AB30").Value = '19:40'

# Row 31
$ws.Range("A31").Value = 112013698
$ws.Range("B31").Value = 77636
$ws.Range("E31").Value = 6425
$ws.Range("F31").Value = 'Garnlav'
$ws.Range("G31").Value = 'Alectoria sarmentosa'
$ws.Range("H31").Value = '(Ach.) Ach.'
$ws.Range("Q31").Value = 610094
$ws.Range("R31").Value = 7121456
$ws.Range("Z31").Value = '19:49'
$ws.Range("AB31").Value = '19:49'

# Row 32
$ws.Range("A32").Value = 112013691
$ws.Range("B32").Value = 88623
$ws.Range("E32").Value = 1962
$ws.Range("F32").Value = 'Vaddporing'
$ws.Range("G32").Value = 'Anomoporia kamtschatica'
$ws.Range("H32").Value = '(Parmasto) Bondartseva'
$ws.Range("Q32").Value = 610134
$ws.Range("R32").Value = 7121461
$ws.Range("Z32").Value = '19:29'
$ws.Range("AB32").Value = '19:29'

# Row 33
$ws.Range("A33").Value = 112013699
$ws.Range("B33").Value = 77636
$ws.Range("E33").Value = 6425
$ws.Range("F33").Value = 'Garnlav'
$ws.Range("G33").Value = 'Alectoria sarmentosa'
$ws.Range("H33").Value = '(Ach.) Ach.'
$ws.Range("Q33").Value = 610068
$ws.Range("R33").Value = 7121408
$ws.Range("Z33").Value = '19:40'
$ws.Range("AB33").Value = '19:40'
